# Billing cash input. SpotStockTest print
# Move the "Total"/"Value" summary rows from rows 71-72 up to rows 65-66,
# and clear the old rows 71-72. Also shrink the Daily_Report1 defined name
# range from A1:M72 to A1:M66 to match the new, shorter report extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content for row 65 (previously lived on row 71) ---
$ws.Range("C65").Value = "Total"
$ws.Range("D65").Value = 232710.53
$ws.Range("E65").Value = 3008.49
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 10250.99
$ws.Range("H65").Value = 36935
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 180
$ws.Range("L65").Value = 472

# --- New content for row 66 (previously lived on row 72) ---
$ws.Range("C66").Value = "Value"
$ws.Range("D66").Value = 242961.52

# --- Remove the old content from rows 71-72, keeping C71/C72's styling ---
$ws.Range("C71").ClearContents()
$ws.Range("D71:L71").Clear()
$ws.Range("C72").ClearContents()
$ws.Range("D72").Clear()

# --- Shrink the defined name range to the new, shorter report extent ---
$wb.Names("Daily_Report1").RefersTo = "=Sheet1!`$A`$1:`$M`$66"
